$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-10 with the recomputed smoke-deployment optimization results
# Row 2
$ws.Cells.Item(2, 1).Value = "M1"
$ws.Cells.Item(2, 2).Value = "FY1"
$ws.Cells.Item(2, 3).Value = -0.99476
$ws.Cells.Item(2, 4).Value = 0.012181
$ws.Cells.Item(2, 5).Value = -0.101506
$ws.Cells.Item(2, 6).Value = 91
$ws.Cells.Item(2, 7).Value = 8000
$ws.Cells.Item(2, 8).Value = 120
$ws.Cells.Item(2, 9).Value = 800
$ws.Cells.Item(2, 10).Value = 108.259551
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 81.194663
$ws.Cells.Item(2, 13).Value = 10450
$ws.Cells.Item(2, 14).Value = 90
$ws.Cells.Item(2, 15).Value = 1050
$ws.Cells.Item(2, 16).Value = 8000
$ws.Cells.Item(2, 17).Value = 120
$ws.Cells.Item(2, 18).Value = 75
$ws.Cells.Item(2, 19).Value = 108.759551
$ws.Cells.Item(2, 20).Value = 4.611044

# Row 3
$ws.Cells.Item(3, 1).Value = "M1"
$ws.Cells.Item(3, 2).Value = "FY1"
$ws.Cells.Item(3, 3).Value = -0.99476
$ws.Cells.Item(3, 4).Value = 0.012181
$ws.Cells.Item(3, 5).Value = -0.101506
$ws.Cells.Item(3, 6).Value = 91
$ws.Cells.Item(3, 7).Value = 8000
$ws.Cells.Item(3, 8).Value = 120
$ws.Cells.Item(3, 9).Value = 800
$ws.Cells.Item(3, 10).Value = 108.259551
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 92.020618
$ws.Cells.Item(3, 13).Value = 9470
$ws.Cells.Item(3, 14).Value = 102
$ws.Cells.Item(3, 15).Value = 950
$ws.Cells.Item(3, 16).Value = 8000
$ws.Cells.Item(3, 17).Value = 120
$ws.Cells.Item(3, 18).Value = 90
$ws.Cells.Item(3, 19).Value = 108.759551
$ws.Cells.Item(3, 20).Value = 3.47636

# Row 4
$ws.Cells.Item(4, 1).Value = "M2"
$ws.Cells.Item(4, 2).Value = "FY2"
$ws.Cells.Item(4, 3).Value = -0.97372
$ws.Cells.Item(4, 4).Value = -0.192924
$ws.Cells.Item(4, 5).Value = -0.121033
$ws.Cells.Item(4, 6).Value = 98
$ws.Cells.Item(4, 7).Value = 6650
$ws.Cells.Item(4, 8).Value = 340
$ws.Cells.Item(4, 9).Value = 735
$ws.Cells.Item(4, 10).Value = 56.065207
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 39.245645
$ws.Cells.Item(4, 13).Value = 8255
$ws.Cells.Item(4, 14).Value = 658
$ws.Cells.Item(4, 15).Value = 934.5
$ws.Cells.Item(4, 16).Value = 6650
$ws.Cells.Item(4, 17).Value = 340
$ws.Cells.Item(4, 18).Value = 75
$ws.Cells.Item(4, 19).Value = 56.565207
$ws.Cells.Item(4, 20).Value = 5.260544

# Row 5
$ws.Cells.Item(5, 1).Value = "M2"
$ws.Cells.Item(5, 2).Value = "FY2"
$ws.Cells.Item(5, 3).Value = -0.97372
$ws.Cells.Item(5, 4).Value = -0.192924
$ws.Cells.Item(5, 5).Value = -0.121033
$ws.Cells.Item(5, 6).Value = 98
$ws.Cells.Item(5, 7).Value = 6650
$ws.Cells.Item(5, 8).Value = 340
$ws.Cells.Item(5, 9).Value = 735
$ws.Cells.Item(5, 10).Value = 56.065207
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 44.852166
$ws.Cells.Item(5, 13).Value = 7720
$ws.Cells.Item(5, 14).Value = 552
$ws.Cells.Item(5, 15).Value = 868
$ws.Cells.Item(5, 16).Value = 6650
$ws.Cells.Item(5, 17).Value = 340
$ws.Cells.Item(5, 18).Value = 90
$ws.Cells.Item(5, 19).Value = 56.565207
$ws.Cells.Item(5, 20).Value = 4.234165

# Row 6
$ws.Cells.Item(6, 1).Value = "M2"
$ws.Cells.Item(6, 2).Value = "FY2"
$ws.Cells.Item(6, 3).Value = -0.97372
$ws.Cells.Item(6, 4).Value = -0.192924
$ws.Cells.Item(6, 5).Value = -0.121033
$ws.Cells.Item(6, 6).Value = 98
$ws.Cells.Item(6, 7).Value = 6650
$ws.Cells.Item(6, 8).Value = 340
$ws.Cells.Item(6, 9).Value = 735
$ws.Cells.Item(6, 10).Value = 56.065207
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 50.458687
$ws.Cells.Item(6, 13).Value = 7185
$ws.Cells.Item(6, 14).Value = 446
$ws.Cells.Item(6, 15).Value = 801.5
$ws.Cells.Item(6, 16).Value = 6650
$ws.Cells.Item(6, 17).Value = 340
$ws.Cells.Item(6, 18).Value = 105
$ws.Cells.Item(6, 19).Value = 56.565207
$ws.Cells.Item(6, 20).Value = 3.746194

# Row 7
$ws.Cells.Item(7, 1).Value = "M3"
$ws.Cells.Item(7, 2).Value = "FY3"
$ws.Cells.Item(7, 3).Value = -0.198479
$ws.Cells.Item(7, 4).Value = 0.979161
$ws.Cells.Item(7, 5).Value = -0.043004
$ws.Cells.Item(7, 6).Value = 81.10075
$ws.Cells.Item(7, 7).Value = 5400
$ws.Cells.Item(7, 8).Value = -40
$ws.Cells.Item(7, 9).Value = 570
$ws.Cells.Item(7, 10).Value = 37.274566
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 27.955924
$ws.Cells.Item(7, 13).Value = 5550
$ws.Cells.Item(7, 14).Value = -780
$ws.Cells.Item(7, 15).Value = 602.5
$ws.Cells.Item(7, 16).Value = 5400
$ws.Cells.Item(7, 17).Value = -40
$ws.Cells.Item(7, 18).Value = 75
$ws.Cells.Item(7, 19).Value = 37.774566
$ws.Cells.Item(7, 20).Value = 4.895187

# Row 8
$ws.Cells.Item(8, 1).Value = "M3"
$ws.Cells.Item(8, 2).Value = "FY3"
$ws.Cells.Item(8, 3).Value = -0.198479
$ws.Cells.Item(8, 4).Value = 0.979161
$ws.Cells.Item(8, 5).Value = -0.043004
$ws.Cells.Item(8, 6).Value = 81.10075
$ws.Cells.Item(8, 7).Value = 5400
$ws.Cells.Item(8, 8).Value = -40
$ws.Cells.Item(8, 9).Value = 570
$ws.Cells.Item(8, 10).Value = 37.274566
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 31.683381
$ws.Cells.Item(8, 13).Value = 5490
$ws.Cells.Item(8, 14).Value = -484
$ws.Cells.Item(8, 15).Value = 589.5
$ws.Cells.Item(8, 16).Value = 5400
$ws.Cells.Item(8, 17).Value = -40
$ws.Cells.Item(8, 18).Value = 90
$ws.Cells.Item(8, 19).Value = 37.774566
$ws.Cells.Item(8, 20).Value = 4.51467

# Row 9
$ws.Cells.Item(9, 1).Value = "M1"
$ws.Cells.Item(9, 2).Value = "FY4"
$ws.Cells.Item(9, 3).Value = -0.93576
$ws.Cells.Item(9, 4).Value = -0.288526
$ws.Cells.Item(9, 5).Value = -0.202748
$ws.Cells.Item(9, 6).Value = 112
$ws.Cells.Item(9, 7).Value = 5000
$ws.Cells.Item(9, 8).Value = 150
$ws.Cells.Item(9, 9).Value = 500
$ws.Cells.Item(9, 10).Value = 57.249134
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 45.799307
$ws.Cells.Item(9, 13).Value = 6200
$ws.Cells.Item(9, 14).Value = 520
$ws.Cells.Item(9, 15).Value = 760
$ws.Cells.Item(9, 16).Value = 5000
$ws.Cells.Item(9, 17).Value = 150
$ws.Cells.Item(9, 18).Value = 75
$ws.Cells.Item(9, 19).Value = 57.749134
$ws.Cells.Item(9, 20).Value = 5.204866

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "FY5"
$ws.Cells.Item(10, 3).Value = -0.96648
$ws.Cells.Item(10, 4).Value = 0.239519
$ws.Cells.Item(10, 5).Value = -0.092446
$ws.Cells.Item(10, 6).Value = 119
$ws.Cells.Item(10, 7).Value = 3800
$ws.Cells.Item(10, 8).Value = 280
$ws.Cells.Item(10, 9).Value = 420
$ws.Cells.Item(10, 10).Value = 79.992232
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 63.993785
$ws.Cells.Item(10, 13).Value = 5640
$ws.Cells.Item(10, 14).Value = -176
$ws.Cells.Item(10, 15).Value = 596
$ws.Cells.Item(10, 16).Value = 3800
$ws.Cells.Item(10, 17).Value = 280
$ws.Cells.Item(10, 18).Value = 75
$ws.Cells.Item(10, 19).Value = 80.492232
$ws.Cells.Item(10, 20).Value = 4.305364

# Remove now-unused extra rows (former rows 11-14) so data ends at row 10
$ws.Range("A11:A14").EntireRow.Delete()
